# Updates market-derived price/profit figures on several Leve sheets.
# Values below were recomputed from a refreshed market-data pull (scheduled
# runner); only numeric result columns (H-N) change, source leve columns
# (A-G) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 301.5
$ws.Range("I33").Value = 109.82609
$ws.Range("K33").Value = 109.82609
$ws.Range("M33").Value = 119.17391

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3244.8438
$ws.Range("I61").Value = 1195
$ws.Range("K61").Value = 1195
$ws.Range("M61").Value = -983

$ws.Range("H136").Value = 3244.8438
$ws.Range("I136").Value = 1195
$ws.Range("K136").Value = 3585
$ws.Range("M136").Value = -1035

$ws.Range("H139").Value = 26736.8
$ws.Range("J139").Value = 26736.8
$ws.Range("L139").Value = 26736.8
$ws.Range("N139").Value = -37016.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 2784.6667
$ws.Range("J17").Value = 2784.6667
$ws.Range("L17").Value = 2784.6667
$ws.Range("N17").Value = -3128.6667

$ws.Range("H81").Value = 8163.3335
$ws.Range("J81").Value = 8163.3335
$ws.Range("L81").Value = 8163.3335
$ws.Range("N81").Value = -10285.3335

$ws.Range("H84").Value = 8163.3335
$ws.Range("J84").Value = 8163.3335
$ws.Range("L84").Value = 24490.0005
$ws.Range("N84").Value = -35098.00049999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 19875.812
$ws.Range("I2").Value = 2999.3333
$ws.Range("J2").Value = 23770.385
$ws.Range("K2").Value = 2999.3333
$ws.Range("L2").Value = 23770.385
$ws.Range("M2").Value = -2886.3333
$ws.Range("N2").Value = -23996.385

$ws.Range("H3").Value = 47858
$ws.Range("I3").Value = 26000
$ws.Range("J3").Value = 56601.2
$ws.Range("K3").Value = 26000
$ws.Range("L3").Value = 56601.2
$ws.Range("M3").Value = -25887
$ws.Range("N3").Value = -56827.2

$ws.Range("H41").Value = 7618
$ws.Range("J41").Value = 10177
$ws.Range("L41").Value = 10177
$ws.Range("N41").Value = -11033

$ws.Range("H50").Value = 17764
$ws.Range("J50").Value = 17764
$ws.Range("L50").Value = 17764
$ws.Range("N50").Value = -19014

$ws.Range("H51").Value = 8474.75
$ws.Range("J51").Value = 8474.75
$ws.Range("L51").Value = 8474.75
$ws.Range("N51").Value = -9946.75

$ws.Range("H59").Value = 18025.666
$ws.Range("J59").Value = 18025.666
$ws.Range("L59").Value = 18025.666
$ws.Range("N59").Value = -20315.666

$ws.Range("H60").Value = 15800.75
$ws.Range("J60").Value = 15734.333
$ws.Range("L60").Value = 15734.333
$ws.Range("N60").Value = -16756.333

$ws.Range("H61").Value = 8474.75
$ws.Range("J61").Value = 8474.75
$ws.Range("L61").Value = 8474.75
$ws.Range("N61").Value = -9170.75

$ws.Range("H68").Value = 34398.57
$ws.Range("J68").Value = 34398.57
$ws.Range("L68").Value = 34398.57
$ws.Range("N68").Value = -35896.57

$ws.Range("H71").Value = 34398.57
$ws.Range("J71").Value = 34398.57
$ws.Range("L71").Value = 103195.71
$ws.Range("N71").Value = -110683.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1900
$ws.Range("I13").Value = 1500
$ws.Range("J13").Value = 2300
$ws.Range("K13").Value = 4500
$ws.Range("L13").Value = 6900
$ws.Range("M13").Value = -4332
$ws.Range("N13").Value = -7236

$ws.Range("H16").Value = 2150.25
$ws.Range("I16").Value = 250
$ws.Range("J16").Value = 2783.6667
$ws.Range("K16").Value = 750
$ws.Range("L16").Value = 8351.000100000001
$ws.Range("M16").Value = -577
$ws.Range("N16").Value = -8697.000100000001

$ws.Range("H21").Value = 1758.3334
$ws.Range("I21").Value = 50
$ws.Range("J21").Value = 2100
$ws.Range("K21").Value = 150
$ws.Range("L21").Value = 6300
$ws.Range("M21").Value = 23
$ws.Range("N21").Value = -6646

$ws.Range("H24").Value = 1954.9546
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1954.9546
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 5864.8638
$ws.Range("N24").Value = -6324.8638
$ws.Range("M24").ClearContents()

$ws.Range("H25").Value = 3000
$ws.Range("I25").Value = 3000
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = -8831
$ws.Range("N25").Value = -9338

$ws.Range("H26").Value = 16965
$ws.Range("I26").Value = 65
$ws.Range("J26").Value = 25415
$ws.Range("K26").Value = 195
$ws.Range("L26").Value = 76245
$ws.Range("M26").Value = 93
$ws.Range("N26").Value = -76821

$ws.Range("H30").Value = 3000
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 3000
$ws.Range("K30").Value = 9000
$ws.Range("L30").Value = 9000
$ws.Range("M30").Value = -8898
$ws.Range("N30").Value = -9204

$ws.Range("H41").Value = 3001.2
$ws.Range("J41").Value = 3001.2
$ws.Range("L41").Value = 9003.599999999999
$ws.Range("N41").Value = -9679.599999999999

$ws.Range("H49").Value = 2227.1428
$ws.Range("I49").Value = 345
$ws.Range("J49").Value = 2980
$ws.Range("K49").Value = 1035
$ws.Range("L49").Value = 8940
$ws.Range("M49").Value = -879
$ws.Range("N49").Value = -9252

$ws.Range("H70").Value = 3661.2856
$ws.Range("I70").Value = 1869.3334
$ws.Range("K70").Value = 5608.0002
$ws.Range("M70").Value = -5293.0002

$ws.Range("H73").Value = 3661.2856
$ws.Range("I73").Value = 1869.3334
$ws.Range("K73").Value = 5608.0002
$ws.Range("M73").Value = -4516.0002

$ws.Range("H100").Value = 2842.2
$ws.Range("J100").Value = 2842.2
$ws.Range("L100").Value = 8526.599999999999
$ws.Range("N100").Value = -10148.6

$ws.Range("H120").Value = 17121.111
$ws.Range("I120").Value = 13522.5
$ws.Range("K120").Value = 40567.5
$ws.Range("M120").Value = -35729.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 60003
$ws.Range("I19").Value = 50000
$ws.Range("J19").Value = 63337.332
$ws.Range("K19").Value = 50000
$ws.Range("L19").Value = 63337.332
$ws.Range("M19").Value = -49712
$ws.Range("N19").Value = -63913.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 882.5
$ws.Range("I16").Value = 955.4375
$ws.Range("K16").Value = 955.4375
$ws.Range("M16").Value = -785.4375

$ws.Range("H46").Value = 1290
$ws.Range("I46").Value = 420
$ws.Range("J46").Value = 2160
$ws.Range("K46").Value = 420
$ws.Range("L46").Value = 2160
$ws.Range("M46").Value = -232
$ws.Range("N46").Value = -2536

$ws.Range("H136").Value = 1562.6976
$ws.Range("I136").Value = 1015.84
$ws.Range("J136").Value = 2322.2222
$ws.Range("K136").Value = 3047.52
$ws.Range("L136").Value = 6966.6666
$ws.Range("M136").Value = -497.52
$ws.Range("N136").Value = -12066.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2512250
$ws.Range("I3").Value = 10000000
$ws.Range("J3").Value = 16333.333
$ws.Range("K3").Value = 10000000
$ws.Range("L3").Value = 16333.333
$ws.Range("M3").Value = -9999886
$ws.Range("N3").Value = -16561.333
